# locode example - apply update
# Changes:
#  - A9: change indicator "+ " -> "*"
#  - A11: new change indicator "+"
#  - C14/D14: "Aberdeen " -> "Abeerdeen "
#  - F14: "--3----- " -> "1-3--6-- "
#  - Active selection moves to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "*"
$ws.Range("A11").Value = "+"

$ws.Range("C14").Value = "Abeerdeen "
$ws.Range("D14").Value = "Abeerdeen "
$ws.Range("F14").Value = "1-3--6-- "

$ws.Range("C9").Select()
